$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.557.63'
$ws.Range('E2').Value = '  +2.27%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.679.09'
$ws.Range('E3').Value = '  +2.86%  '
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '219.59'
$ws.Range('E5').Value = '  +2.27%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.526'
$ws.Range('E6').Value = '  +1.77%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '29.94'
$ws.Range('E8').Value = '  +4.66%  '
$ws.Range('E9').Value = '  +2.10%  '
$ws.Range('E10').Value = '  +3.06%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0907'
$ws.Range('E11').Value = '  -0.67%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.921.24'
$ws.Range('E12').Value = '  +2.84%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.43'
$ws.Range('E13').Value = '  +12.84%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.620'
$ws.Range('E14').Value = '  +9.75%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '1.674.44'
$ws.Range('E15').Value = '  +2.57%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.97'
$ws.Range('E16').Value = '  +2.70%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.575.16'
$ws.Range('E17').Value = '  +2.26%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '66.32'
$ws.Range('E18').Value = '  +3.22%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '244.54'
$ws.Range('E19').Value = '  +1.12%  '
$ws.Range('E20').Value = '  +2.35%  '
$ws.Range('E21').Value = '  -0.08%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.28'
$ws.Range('E22').Value = '  +3.68%  '
$ws.Range('B23').Value = 'Avalanche'
$ws.Range('C23').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.15'
$ws.Range('E23').Value = '  +2.79%  '
$ws.Range('E24').Value = '  +0.66%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '157.59'
$ws.Range('E25').Value = '  +0.01%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.87'
$ws.Range('E26').Value = '  +2.10%  '
$ws.Range('E27').Value = '  +2.13%  '
$ws.Range('E28').Value = '  +1.50%  '
$ws.Range('E29').Value = '  -0.14%  '
$ws.Range('E30').Value = '  +2.02%  '
$ws.Range('E31').Value = '  +1.85%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.47'
$ws.Range('E32').Value = '  +3.24%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.505.56'
$ws.Range('E33').Value = '  +5.21%  '
$ws.Range('E34').Value = '  +3.75%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.76'
$ws.Range('E35').Value = '  +7.04%  '
$ws.Range('E36').Value = '  -1.08%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '83.69'
$ws.Range('E37').Value = '  +10.33%  '
$ws.Range('E38').Value = '  +4.99%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.589'
$ws.Range('E39').Value = '  +6.63%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.69'
$ws.Range('E40').Value = '  -4.85%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.30'
$ws.Range('E41').Value = '  +0.00%  '
$ws.Range('E42').Value = '  +1.60%  '
$ws.Range('B43').Value = 'Kaspa'
$ws.Range('C43').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0501'
$ws.Range('E43').Value = '  +1.77%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.99'
$ws.Range('E44').Value = '  +0.05%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.03'
$ws.Range('E45').Value = '  +1.12%  '
$ws.Range('E46').Value = '  -0.02%  '
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '5.52'
$ws.Range('E47').Value = '  +3.10%  '
$ws.Range('B48').Value = 'BitcoinSV'
$ws.Range('C48').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '51.25'
$ws.Range('E48').Value = '  -4.41%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.814.20'
$ws.Range('E49').Value = '  +2.10%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '94.52'
$ws.Range('E50').Value = '  +6.12%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0₆0112'
$ws.Range('E51').Value = '  +2.10%  '
